
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Caso = MTOW @ Nivel do Mar)
$ws.Range("B2").Value = 434376.4019098243
$ws.Range("D2").Value = 4117.311866443833
$ws.Range("E2").Value = 71.86433215499773
$ws.Range("F2").Value = 113.6274860683347
$ws.Range("I2").Value = 26.39996022013026
$ws.Range("J2").Value = 0.7328704778298667
$ws.Range("K2").Value = 3.031545613680756
$ws.Range("L2").Value = 1.894716008550473

# Row 3 (Caso = MZFW @ Nivel do Mar)
$ws.Range("B3").Value = 351110.32440492
$ws.Range("D3").Value = 3328.059946966066
$ws.Range("E3").Value = 64.61034829962428
$ws.Range("F3").Value = 102.1579305218
$ws.Range("I3").Value = 21.33932358298599
$ws.Range("J3").Value = 0.7049204794757324
$ws.Range("K3").Value = 3.607444275370693
$ws.Range("L3").Value = 2.254652672106683

# Row 4 (Caso = Peso de Cruzeiro @ 35000 ft)
$ws.Range("D4").Value = 3633.084682464455
$ws.Range("E4").Value = 121.1324471307076
$ws.Range("F4").Value = 191.5272157414821
$ws.Range("I4").Value = 75.0062523101635
$ws.Range("J4").Value = 0.8219223302565979
$ws.Range("K4").Value = 1.19666741498567
$ws.Range("L4").Value = 0.7479171343660437

# Row 5 (Caso = Peso de Pouso @ 10000 ft)
$ws.Range("D5").Value = 3633.084682464455
$ws.Range("E5").Value = 78.54932705867407
$ws.Range("F5").Value = 124.1973910894523
$ws.Range("I5").Value = 31.53998060994971
$ws.Range("J5").Value = 0.7533984132787424
$ws.Range("K5").Value = 2.608575622681944
$ws.Range("L5").Value = 1.630359764176216
